# Update the release/version string throughout the workbook.
#
# Old version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# New version string: "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet ---
$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = "Version: " + $newVersion

$aboutWs.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for KWK So" + [char]347 + "nica Coal Mine, Poland, M1471, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$boundariesWs = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 14; $row++) {
    $boundariesWs.Range("S" + $row).Value = $newVersion
}
